$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.033.30'
$ws.Cells.Item(2, 5).Value = '  -0.61%  '
$ws.Cells.Item(3, 4).Value = '2.217.45'
$ws.Cells.Item(3, 5).Value = '  -1.36%  '
$ws.Cells.Item(4, 5).Value = '  +0.25%  '
$c = $ws.Cells.Item(5, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '242.25'
$c.Style = $origStyle
$ws.Cells.Item(5, 5).Value = '  -2.09%  '
$c = $ws.Cells.Item(6, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.626'
$c.Style = $origStyle
$ws.Cells.Item(6, 5).Value = '  -0.44%  '
$c = $ws.Cells.Item(7, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '73.72'
$c.Style = $origStyle
$ws.Cells.Item(7, 5).Value = '  -1.35%  '
$ws.Cells.Item(8, 5).Value = '  +0.13%  '
$c = $ws.Cells.Item(9, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.612'
$c.Style = $origStyle
$ws.Cells.Item(9, 5).Value = '  -1.21%  '
$c = $ws.Cells.Item(10, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '43.77'
$c.Style = $origStyle
$ws.Cells.Item(10, 5).Value = '  +3.22%  '
$c = $ws.Cells.Item(11, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0960'
$c.Style = $origStyle
$ws.Cells.Item(11, 5).Value = '  +1.49%  '
$c = $ws.Cells.Item(12, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '7.12'
$c.Style = $origStyle
$ws.Cells.Item(12, 5).Value = '  -1.14%  '
$c = $ws.Cells.Item(13, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.103'
$c.Style = $origStyle
$ws.Cells.Item(13, 5).Value = '  +0.48%  '
$ws.Cells.Item(14, 4).Value = '2.555.64'
$ws.Cells.Item(14, 5).Value = '  -1.13%  '
$c = $ws.Cells.Item(15, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '14.27'
$c.Style = $origStyle
$ws.Cells.Item(15, 5).Value = '  -1.96%  '
$c = $ws.Cells.Item(16, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.845'
$c.Style = $origStyle
$ws.Cells.Item(16, 5).Value = '  -1.50%  '
$ws.Cells.Item(17, 4).Value = '2.233.50'
$ws.Cells.Item(17, 5).Value = '  -0.62%  '
$ws.Cells.Item(18, 4).Value = '41.892.48'
$ws.Cells.Item(18, 5).Value = '  -0.54%  '
$c = $ws.Cells.Item(19, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0000109'
$c.Style = $origStyle
$ws.Cells.Item(19, 5).Value = '  +11.02%  '
$ws.Cells.Item(20, 2).Value = 'Litecoin'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Cells.Item(20, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '72.64'
$c.Style = $origStyle
$ws.Cells.Item(20, 5).Value = '  +0.89%  '
$ws.Cells.Item(21, 2).Value = 'Uniswap'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Cells.Item(21, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '6.16'
$c.Style = $origStyle
$ws.Cells.Item(21, 5).Value = '  +0.07%  '
$c = $ws.Cells.Item(22, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '10.55'
$c.Style = $origStyle
$ws.Cells.Item(22, 5).Value = '  +18.68%  '
$c = $ws.Cells.Item(23, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '229.66'
$c.Style = $origStyle
$ws.Cells.Item(23, 5).Value = '  -1.10%  '
$c = $ws.Cells.Item(24, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.09'
$c.Style = $origStyle
$ws.Cells.Item(24, 5).Value = '  -6.99%  '
$c = $ws.Cells.Item(25, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '11.56'
$c.Style = $origStyle
$ws.Cells.Item(25, 5).Value = '  +2.82%  '
$ws.Cells.Item(26, 5).Value = '  -0.05%  '
$ws.Cells.Item(27, 5).Value = '  -1.30%  '
$c = $ws.Cells.Item(28, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.26'
$c.Style = $origStyle
$ws.Cells.Item(28, 5).Value = '  -2.29%  '
$ws.Cells.Item(29, 5).Value = '  +2.50%  '
$c = $ws.Cells.Item(30, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '166.89'
$c.Style = $origStyle
$ws.Cells.Item(30, 5).Value = '  -1.50%  '
$c = $ws.Cells.Item(31, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '20.60'
$c.Style = $origStyle
$ws.Cells.Item(31, 5).Value = '  -0.42%  '
$c = $ws.Cells.Item(32, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.63'
$c.Style = $origStyle
$ws.Cells.Item(32, 5).Value = '  +7.28%  '
$c = $ws.Cells.Item(33, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0796'
$c.Style = $origStyle
$ws.Cells.Item(33, 5).Value = '  -4.12%  '
$ws.Cells.Item(34, 5).Value = '  +0.14%  '
$c = $ws.Cells.Item(35, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.113'
$c.Style = $origStyle
$ws.Cells.Item(35, 5).Value = '  -5.37%  '
$c = $ws.Cells.Item(36, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '28.87'
$c.Style = $origStyle
$ws.Cells.Item(36, 5).Value = '  -4.71%  '
$c = $ws.Cells.Item(37, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '4.27'
$c.Style = $origStyle
$ws.Cells.Item(37, 5).Value = '  -5.79%  '
$c = $ws.Cells.Item(38, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.0303'
$c.Style = $origStyle
$ws.Cells.Item(38, 5).Value = '  -0.58%  '
$c = $ws.Cells.Item(39, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '12.87'
$c.Style = $origStyle
$ws.Cells.Item(39, 5).Value = '  -4.78%  '
$c = $ws.Cells.Item(40, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '66.54'
$c.Style = $origStyle
$ws.Cells.Item(40, 5).Value = '  +7.02%  '
$c = $ws.Cells.Item(41, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.12'
$c.Style = $origStyle
$ws.Cells.Item(41, 5).Value = '  -3.24%  '
$c = $ws.Cells.Item(42, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '5.64'
$c.Style = $origStyle
$ws.Cells.Item(42, 5).Value = '  -2.74%  '
$c = $ws.Cells.Item(43, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '0.202'
$c.Style = $origStyle
$ws.Cells.Item(43, 5).Value = '  -0.37%  '
$c = $ws.Cells.Item(44, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '8.76'
$c.Style = $origStyle
$ws.Cells.Item(44, 5).Value = '  +1.31%  '
$c = $ws.Cells.Item(45, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '104.08'
$c.Style = $origStyle
$ws.Cells.Item(45, 5).Value = '  -4.94%  '
$ws.Cells.Item(46, 5).Value = '  -0.51%  '
$c = $ws.Cells.Item(47, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '2.45'
$c.Style = $origStyle
$ws.Cells.Item(47, 5).Value = '  +5.83%  '
$ws.Cells.Item(48, 5).Value = '  -0.48%  '
$c = $ws.Cells.Item(49, 4)
$origStyle = $c.Style
$c.NumberFormat = '@'
$c.Value = '1.16'
$c.Style = $origStyle
$ws.Cells.Item(49, 5).Value = '  -0.16%  '
$ws.Cells.Item(50, 5).Value = '  +0.22%  '
$ws.Cells.Item(51, 4).Value = '2.427.62'
$ws.Cells.Item(51, 5).Value = '  -1.25%  '
